$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading, then work on the bullet
# paragraphs that immediately follow its "Impact" sub-heading. Note: several
# of these bullet strings (e.g. the 87% turnout-accuracy one) also appear
# verbatim earlier in the "PROFESSIONAL EXPERIENCE" section, so we must not
# do a blind document-wide find/replace — only the bullets inside this
# specific section should change.

$count = $d.Paragraphs.Count
$sectionStart = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $sectionStart = $i
        break
    }
}
if ($sectionStart -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# The first bullet paragraph after the heading + its "Impact" sub-heading.
$firstBullet = $sectionStart + 2
if ($d.Paragraphs.Item($firstBullet).Range.Text -notmatch "Discovered systematic race coding errors") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 1: $($d.Paragraphs.Item($firstBullet).Range.Text)"
}

$p1 = $d.Paragraphs.Item($firstBullet)
$p1.Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"

$p2 = $d.Paragraphs.Item($firstBullet + 1)
if ($p2.Range.Text -notmatch "Achieved 87% prediction accuracy") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 2: $($p2.Range.Text)"
}
$p2.Range.Text = "• Reduced polling margins from ±4.2% to ±2.1%"

$p3 = $d.Paragraphs.Item($firstBullet + 2)
if ($p3.Range.Text -notmatch "Built redistricting platform") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 3: $($p3.Range.Text)"
}
$p3.Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"

$p4 = $d.Paragraphs.Item($firstBullet + 3)
if ($p4.Range.Text -notmatch "Developed longitudinal data analysis methods") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 4: $($p4.Range.Text)"
}
$p4.Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# Remove the two trailing bullets (expert testimony / FEC compliance)
# entirely — they're job duties, not accomplishments, and the section is
# being trimmed to four punchy, single-line bullets.
$p5 = $d.Paragraphs.Item($firstBullet + 4)
if ($p5.Range.Text -notmatch "Provided expert testimony") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 5: $($p5.Range.Text)"
}
$p5.Range.Delete()

$p6 = $d.Paragraphs.Item($firstBullet + 4)
if ($p6.Range.Text -notmatch "Demystified FEC compliance") {
    throw "Unexpected paragraph at KEY ACHIEVEMENTS bullet 6: $($p6.Range.Text)"
}
$p6.Range.Delete()

Write-Output "Done."
